# Update cryptocurrency price/volume data per Feb 21 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.698.42"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "2.897.18"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Formula = "'367.18"
$ws.Range("E5").Value = "  +4.83%  "
$ws.Range("D6").Formula = "'101.89"
$ws.Range("E6").Value = "  -3.98%  "
$ws.Range("D7").Formula = "'0.538"
$ws.Range("E7").Value = "  -3.13%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").Formula = "'36.24"
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Formula = "'0.0828"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Formula = "'18.18"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").Value = "3.349.54"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").Value = "2.893.52"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Formula = "'0.917"
$ws.Range("E17").Value = "  -4.74%  "
$ws.Range("D18").Value = "50.702.46"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("E19").Value = "  -5.89%  "
$ws.Range("D20").Formula = "'7.13"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("D21").Formula = "'12.78"
$ws.Range("E21").Value = "  -4.43%  "
$ws.Range("D22").Value = "0.0₃0937"
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("D23").Formula = "'67.82"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").Formula = "'256.82"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Formula = "'4.23"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Formula = "'0.166"
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("D29").Formula = "'25.39"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("D30").Formula = "'7.02"
$ws.Range("E30").Value = "  -3.31%  "
$ws.Range("D31").Formula = "'0.101"
$ws.Range("E31").Value = "  -5.23%  "
$ws.Range("D32").Formula = "'6.18"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").Formula = "'9.80"
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("D35").Formula = "'50.70"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Formula = "'33.87"
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("E39").Value = "  -5.23%  "
$ws.Range("D40").Formula = "'16.87"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("D41").Formula = "'2.56"
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("D42").Formula = "'1.83"
$ws.Range("E42").Value = "  -5.42%  "
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("D44").Formula = "'118.57"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Formula = "'21.71"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Formula = "'2.08"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").Value = "2.008.28"
$ws.Range("E48").Value = "  -4.11%  "
$ws.Range("D49").Formula = "'3.12"
$ws.Range("E49").Value = "  -5.35%  "
$ws.Range("D50").Value = "3.182.64"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("E51").Value = "  -1.43%  "
